# Update "想去人数" (interest count) values in column F on sheets
# "展览" and "全部类型" to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (row number -> new F value)
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    7  = 4622
    8  = 194
    9  = 125
    15 = 1026
    16 = 82
    21 = 97
    22 = 3607
    29 = 3386
    32 = 2486
    35 = 127
    36 = 226
    40 = 1530
    41 = 917
    42 = 22
    45 = 76
    48 = 554
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# Sheet "全部类型" (row number -> new F value)
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    7  = 4622
    8  = 194
    9  = 125
    16 = 1026
    17 = 82
    22 = 97
    23 = 3607
    30 = 3386
    33 = 2486
    36 = 127
    37 = 226
    41 = 1530
    42 = 917
    43 = 22
    46 = 78
    49 = 554
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
